# Atualização de bases das ligas, do dia: 27-04-2024 às 11:27
#
# The source data rows got re-fetched/re-ordered upstream; for a handful
# of match pairs the two rows' data (everything except the running index
# in column A) ended up swapped. Also one row got a small odds refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowData {
    param([int]$Row1, [int]$Row2, [int]$FirstCol, [int]$LastCol)

    for ($c = $FirstCol; $c -le $LastCol; $c++) {
        $cell1 = $ws.Cells.Item($Row1, $c)
        $cell2 = $ws.Cells.Item($Row2, $c)
        $v1 = $cell1.Value()
        $v2 = $cell2.Value()
        $cell1.Value = $v2
        $cell2.Value = $v1
    }
}

# Rows whose entire records (minus col A, the running index) got swapped
# with their neighbour row.
Swap-RowData 38 39 2 28
Swap-RowData 110 111 2 28
Swap-RowData 129 130 2 28
Swap-RowData 237 238 2 28

# Row 250: closing-line odds refresh (no row swap here).
$ws.Range("Q250").Value = 2.025
$ws.Range("R250").Value = 1.775
$ws.Range("T250").Value = 1.95
$ws.Range("U250").Value = 1.85
